# Updated Forecast with Timespan
#
# For each forecast sheet (one per market segment), row 24 (Timeline = 42704)
# gets new "Room Nights Sold" / "Average Room Rate" / "Revenue" actuals in
# columns B/H/N, which are then mirrored into the adjoining
# Forecast/Lower-CB/Upper-CB columns (C-E / I-K / O-Q) for that same row.
#
# Row 25 (Timeline = 42735) changes from holding a literal actual value in
# B/H/N to being blank there (those become forecast-only rows), while its
# Forecast/Lower-CB/Upper-CB columns (C-E / I-K / O-Q) are normalized to 2686.

$wb = $excel.ActiveWorkbook

function Set-ForecastRow24 {
    param(
        [string]$SheetName,
        [double]$B24,
        [double]$H24,
        [double]$N24
    )

    $ws = $wb.Worksheets.Item($SheetName)

    $ws.Range("B24").Value = $B24
    $ws.Range("C24").Value = $B24
    $ws.Range("D24").Value = $B24
    $ws.Range("E24").Value = $B24

    $ws.Range("H24").Value = $H24
    $ws.Range("I24").Value = $H24
    $ws.Range("J24").Value = $H24
    $ws.Range("K24").Value = $H24

    $ws.Range("N24").Value = $N24
    $ws.Range("O24").Value = $N24
    $ws.Range("P24").Value = $N24
    $ws.Range("Q24").Value = $N24
}

function Set-ForecastRow25 {
    param(
        [string]$SheetName
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # B25 / H25 / N25 no longer hold a literal actual - clear them out.
    $ws.Range("B25").ClearContents()
    $ws.Range("H25").ClearContents()
    $ws.Range("N25").ClearContents()

    # Forecast / Lower-CB / Upper-CB columns are normalized to 2686.
    $ws.Range("C25").Value = 2686
    $ws.Range("D25").Value = 2686
    $ws.Range("E25").Value = 2686

    $ws.Range("I25").Value = 2686
    $ws.Range("J25").Value = 2686
    $ws.Range("K25").Value = 2686

    $ws.Range("O25").Value = 2686
    $ws.Range("P25").Value = 2686
    $ws.Range("Q25").Value = 2686
}

# Sheet name -> (new B24, new H24, new N24)
# (named parameters aren't reliably bound by this host's parser, so these
# calls use positional arguments)
Set-ForecastRow24 "RCK"       514   10693 10559
Set-ForecastRow24 "CORP"      566   5974  31381
Set-ForecastRow24 "CORPO"     531   16667 917
Set-ForecastRow24 "INDO"      0     0     0
Set-ForecastRow24 "INDR"      0     0     0
Set-ForecastRow24 "PKG-PRM"   17583 0     0
Set-ForecastRow24 "WSOL"      0     0     0
Set-ForecastRow24 "WSOF"      0     0     0
Set-ForecastRow24 "CON-ASSOC" 0     0     0
Set-ForecastRow24 "CORPM"     0     0     0
Set-ForecastRow24 "GRPO"      0     0     0
Set-ForecastRow24 "GRPT"      0     0     0

Set-ForecastRow25 "RCK"
Set-ForecastRow25 "CORP"
Set-ForecastRow25 "CORPO"
Set-ForecastRow25 "INDO"
Set-ForecastRow25 "INDR"
Set-ForecastRow25 "PKG-PRM"
Set-ForecastRow25 "WSOL"
Set-ForecastRow25 "WSOF"
Set-ForecastRow25 "CON-ASSOC"
Set-ForecastRow25 "CORPM"
Set-ForecastRow25 "GRPO"
Set-ForecastRow25 "GRPT"

# QD and GOV-NGO sheets are intentionally left untouched - they already
# reflect the target layout (blank B/H/N25, 2686 in the forecast columns).
